# Update the dSF column (column F) values to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = -2
    4  = -1
    5  = -1
    7  = -3
    8  = -2
    9  = 1
    10 = 1
    11 = -1
    12 = -5
    13 = -2
    14 = -2
    15 = -3
    18 = -4
    20 = -3
    21 = 3
    22 = 3
    23 = -2
    24 = 1
    25 = 1
    26 = 4
    27 = 1
    28 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
